# Regenerate handback-status report:
#  - refresh timestamps/hash for existing file ab856b28... -> 1aa29009...
#  - add a brand new handed-back file 603718cb-1111-4a69-ba0a-989b0d347a7d.md
#    as a new row (row 3) on every sheet (Overview, zh-cn, de-de)

$wb = $excel.ActiveWorkbook

$uuid1 = "1aa29009-39e0-4b33-a645-3f348e20e891"
$uuid2 = "603718cb-1111-4a69-ba0a-989b0d347a7d"
$hash1 = "1d419a78037f0c5f01dfa176c821250c8473c753"
$hash2 = "57328d7613f7bf05c785a2af73361c52d54d9c34"

# Leading apostrophe forces Excel to store the value as literal text instead
# of auto-converting to a boolean / number. A lone apostrophe yields an
# empty text cell (instead of clearing the cell entirely).
$txtTrue  = "'True"
$txtFalse = "'False"
$txtEmpty = "'"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Update existing row 2 values (file renamed + regenerated date)
$wsOverview.Range("A2").Value = "$uuid1.md"
$wsOverview.Range("B2").Value = "e2e\$uuid1.md"
$wsOverview.Range("G2").Value = "2016-08-26 17:02:23"

# Add a new table row for the new file
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "$uuid2.md"
$wsOverview.Range("B3").Value = "e2e\$uuid2.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").NumberFormat = $dateFmt
$wsOverview.Range("G3").Value = "2016-08-26 17:02:23"

# Refresh hyperlinks (old uuid1 link + brand new uuid2 link)
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6488b6cc3317a7f4893997d6b4a8b7660b32ca81/e2e/$uuid1.md", $null, $null, "e2e\$uuid1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6488b6cc3317a7f4893997d6b4a8b7660b32ca81/e2e/$uuid2.md", $null, $null, "e2e\$uuid2.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Update existing row 2 (file renamed, regenerated hash + timestamps for uuid1)
$wsZhCn.Range("A2").Value = "$uuid1.md"
$wsZhCn.Range("G2").Value = "$uuid1.$hash1.zh-cn.xlf"
$wsZhCn.Range("H2").NumberFormat = $dateFmt
$wsZhCn.Range("H2").Value = "2016-08-26 17:02:19"
$wsZhCn.Range("I2").Value = "$uuid1.md"
$wsZhCn.Range("J2").Value = "$uuid1.$hash1.zh-cn.xlf"
$wsZhCn.Range("K2").NumberFormat = $dateFmt
$wsZhCn.Range("K2").Value = "2016-08-26 17:02:36"

# Add a new table row for the new file
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = "$uuid2.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = $txtTrue
$wsZhCn.Range("G3").Value = "$uuid2.$hash2.zh-cn.xlf"
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("H3").Value = "2016-08-26 17:02:19"
$wsZhCn.Range("I3").Value = "$uuid2.md"
$wsZhCn.Range("J3").Value = "$uuid2.$hash2.zh-cn.xlf"
$wsZhCn.Range("K3").NumberFormat = $dateFmt
$wsZhCn.Range("K3").Value = "2016-08-26 17:02:36"
$wsZhCn.Range("L3").Value = $txtEmpty
$wsZhCn.Range("M3").Value = $txtTrue
$wsZhCn.Range("N3").Value = $txtEmpty
$wsZhCn.Range("O3").Value = $txtFalse
$wsZhCn.Range("P3").Value = $txtEmpty

# Refresh hyperlinks (old uuid1 links + brand new uuid2 links)
$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Range("I2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6488b6cc3317a7f4893997d6b4a8b7660b32ca81/e2e/$uuid1.md", $null, $null, "$uuid1.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2e8c8f7f31cf377833b823154a9d75c2abc6cd8a/e2e/$uuid1.md", $null, $null, "$uuid1.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6488b6cc3317a7f4893997d6b4a8b7660b32ca81/e2e/$uuid2.md", $null, $null, "$uuid2.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2e8c8f7f31cf377833b823154a9d75c2abc6cd8a/e2e/$uuid2.md", $null, $null, "$uuid2.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update existing row 2 (file renamed, regenerated hash + timestamps for uuid1)
$wsDeDe.Range("A2").Value = "$uuid1.md"
$wsDeDe.Range("G2").Value = "$uuid1.$hash1.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-26 17:02:23"
$wsDeDe.Range("I2").Value = "$uuid1.md"
$wsDeDe.Range("J2").Value = "$uuid1.$hash1.de-de.xlf"
$wsDeDe.Range("K2").NumberFormat = $dateFmt
$wsDeDe.Range("K2").Value = "2016-08-26 17:02:43"

# Add a new table row for the new file
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = "$uuid2.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = $txtTrue
$wsDeDe.Range("G3").Value = "$uuid2.$hash2.de-de.xlf"
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("H3").Value = "2016-08-26 17:02:23"
$wsDeDe.Range("I3").Value = "$uuid2.md"
$wsDeDe.Range("J3").Value = "$uuid2.$hash2.de-de.xlf"
$wsDeDe.Range("K3").NumberFormat = $dateFmt
$wsDeDe.Range("K3").Value = "2016-08-26 17:02:43"
$wsDeDe.Range("L3").Value = $txtEmpty
$wsDeDe.Range("M3").Value = $txtTrue
$wsDeDe.Range("N3").Value = $txtEmpty
$wsDeDe.Range("O3").Value = $txtFalse
$wsDeDe.Range("P3").Value = $txtEmpty

# Refresh hyperlinks (old uuid1 links + brand new uuid2 links)
$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Range("I2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6488b6cc3317a7f4893997d6b4a8b7660b32ca81/e2e/$uuid1.md", $null, $null, "$uuid1.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/5957d1abf9c2974ab02a1799bbacbd7507655ebc/e2e/$uuid1.md", $null, $null, "$uuid1.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6488b6cc3317a7f4893997d6b4a8b7660b32ca81/e2e/$uuid2.md", $null, $null, "$uuid2.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/5957d1abf9c2974ab02a1799bbacbd7507655ebc/e2e/$uuid2.md", $null, $null, "$uuid2.md") | Out-Null

Write-Host "Report regenerated."
